# Update "想去人数" (want-to-go count) values in column F across the
# four worksheets, matching the upstream data refresh at commit 456a3b4.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 207  # was 205
$ws.Range("F5").Value = 1381  # was 1380
$ws.Range("F8").Value = 697949  # was 697886
$ws.Range("F9").Value = 1776  # was 1775
$ws.Range("F15").Value = 2422  # was 2418
$ws.Range("F16").Value = 1180  # was 1179
$ws.Range("F17").Value = 2798  # was 2797
$ws.Range("F19").Value = 974  # was 969
$ws.Range("F20").Value = 1629  # was 1628
$ws.Range("F22").Value = 562  # was 561
$ws.Range("F23").Value = 1119  # was 1117
$ws.Range("F24").Value = 1371  # was 1369
$ws.Range("F25").Value = 1131  # was 1130
$ws.Range("F29").Value = 1526  # was 1525
$ws.Range("F31").Value = 1394  # was 1395
$ws.Range("F32").Value = 3680  # was 3676
$ws.Range("F34").Value = 1175  # was 1174
$ws.Range("F37").Value = 208  # was 207
$ws.Range("F38").Value = 133  # was 131
$ws.Range("F42").Value = 3232  # was 3230
$ws.Range("F43").Value = 1048  # was 1047

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 134  # was 132
$ws.Range("F8").Value = 147271  # was 147263
$ws.Range("F9").Value = 147271  # was 147263
$ws.Range("F19").Value = 207  # was 206

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F9").Value = 663  # was 662
$ws.Range("F10").Value = 1638  # was 1637
$ws.Range("F12").Value = 2147  # was 2146

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 663  # was 662
$ws.Range("F5").Value = 1638  # was 1637
$ws.Range("F6").Value = 207  # was 205
$ws.Range("F7").Value = 192  # was 0
$ws.Range("F8").Value = 2147  # was 2146
$ws.Range("F9").Value = 1381  # was 1380
$ws.Range("F12").Value = 697950  # was 697886
$ws.Range("F13").Value = 134  # was 132
$ws.Range("F15").Value = 1776  # was 1775
$ws.Range("F16").Value = 147271  # was 147263
$ws.Range("F21").Value = 2422  # was 2418
$ws.Range("F22").Value = 1180  # was 1179
$ws.Range("F23").Value = 2798  # was 2797
$ws.Range("F25").Value = 974  # was 969
$ws.Range("F27").Value = 1629  # was 1628
$ws.Range("F28").Value = 562  # was 561
$ws.Range("F29").Value = 207  # was 206
$ws.Range("F30").Value = 1119  # was 1117
$ws.Range("F31").Value = 1371  # was 1369
$ws.Range("F32").Value = 1131  # was 1130
$ws.Range("F35").Value = 1526  # was 1525
$ws.Range("F36").Value = 1394  # was 1395
$ws.Range("F37").Value = 3680  # was 3676
$ws.Range("F39").Value = 1175  # was 1174
$ws.Range("F44").Value = 133  # was 131
$ws.Range("F48").Value = 3232  # was 3230
$ws.Range("F50").Value = 1048  # was 1047
